# Update the Bangladesh life-insurance capital-structure database:
# replace rows 2-9 of data with the refreshed dataset (rows 2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous data rows (2-9); row 1 (headers) is left untouched.
$ws.Range("A2:AQ9").ClearContents()

# Row 2: summary row
$ws.Range("A2").Value = 'Bangladesh'
$ws.Range("B2").Value = "'5"
$ws.Range("C2").Value = 'Insurance (Life)'
$ws.Range("D2").Value = -0.03195
$ws.Range("E2").Value = -0.07169999999999999
$ws.Range("G2").Value = 0.1607189119170984
$ws.Range("H2").Value = 0.1607189119170984
$ws.Range("I2").Value = 0.1702072538860104
$ws.Range("J2").Value = 0.1589233645370776
$ws.Range("K2").Value = 43.31
$ws.Range("L2").Value = 0.1402525906735751
$ws.Range("M2").Value = 7.646
$ws.Range("N2").Value = 0.01854744808849214
$ws.Range("O2").Value = 0.1765412145001154
$ws.Range("P2").Value = 7.646
$ws.Range("Q2").Value = 0.01854744808849214
$ws.Range("R2").Value = 0.1765412145001154
$ws.Range("S2").Value = 0.0
$ws.Range("T2").Value = 0.0
$ws.Range("U2").Value = 378.36
$ws.Range("V2").Value = 0.9178148651271104
$ws.Range("W2").Value = 1.043956043956044
$ws.Range("X2").Value = 0.08846408071440458
$ws.Range("Y2").Value = 0.9554919632416392
$ws.Range("Z2").Value = -1.359573812354159
$ws.Range("AA2").Value = -0.07812183383991894
$ws.Range("AB2").Value = 0.08846408071440458
$ws.Range("AC2").Value = -0.1665859145543235
$ws.Range("AD2").Value = 0.096
$ws.Range("AE2").Value = 0.0
$ws.Range("AF2").Value = 0.096
$ws.Range("AG2").Value = -378.264
$ws.Range("AH2").Value = 0.0002328198362500485
$ws.Range("AI2").Value = 0.0006814503535023708
$ws.Range("AJ2").Value = -11.13327054391335
$ws.Range("AK2").Value = 1.592797830590692
$ws.Range("AL2").Value = 0.047
$ws.Range("AM2").Value = 0.047
$ws.Range("AN2").Value = 0.001711229946524064
$ws.Range("AO2").Value = 1118.297872340426
$ws.Range("AP2").Value = -6.742673796791444
$ws.Range("AQ2").Value = 1118.297872340426

# Row 3: Sunlife Insurance Company Limited (DSE:SUNLIFEINS)
$ws.Range("A3").Value = 'Bangladesh'
$ws.Range("B3").Value = 'Sunlife Insurance Company Limited (DSE:SUNLIFEINS)'
$ws.Range("C3").Value = 'Insurance (Life)'
$ws.Range("D3").Value = -0.142
$ws.Range("G3").Value = -0.5
$ws.Range("H3").Value = -0.5
$ws.Range("I3").Value = -0.5578431372549021
$ws.Range("J3").Value = -0.5578431372549021
$ws.Range("K3").Value = -5.81
$ws.Range("L3").Value = -0.5696078431372549
$ws.Range("M3").Value = -0.0
$ws.Range("N3").Value = -0.0
$ws.Range("O3").Value = 0.0
$ws.Range("P3").Value = -0.0
$ws.Range("Q3").Value = -0.0
$ws.Range("R3").Value = 0.0
$ws.Range("S3").Value = 0.0
$ws.Range("U3").Value = 4.86
$ws.Range("V3").Value = 0.4939024390243903
$ws.Range("W3").Value = -1.373522458628841
$ws.Range("X3").Value = 0.08904333008548558
$ws.Range("Y3").Value = -1.462565788714327
$ws.Range("Z3").Value = -8.360655737704919
$ws.Range("AA3").Value = 4.66393442622951
$ws.Range("AB3").Value = 0.08860047459545936
$ws.Range("AC3").Value = 4.57533395163405
$ws.Range("AD3").Value = 0.096
$ws.Range("AE3").Value = 0.0
$ws.Range("AF3").Value = 0.096
$ws.Range("AG3").Value = -4.764
$ws.Range("AH3").Value = 0.00966183574879227
$ws.Range("AI3").Value = 0.02224281742354032
$ws.Range("AJ3").Value = -0.938534278959811
$ws.Range("AK3").Value = 8.757352941176464
$ws.Range("AL3").Value = 0.0
$ws.Range("AM3").Value = 0.0
$ws.Range("AN3").Value = -0.0187866927592955
$ws.Range("AP3").Value = 0.9322896281800391

# Row 4: Sandhani Life Insurance Company Limited (DSE:SANDHANINS)
$ws.Range("A4").Value = 'Bangladesh'
$ws.Range("B4").Value = 'Sandhani Life Insurance Company Limited (DSE:SANDHANINS)'
$ws.Range("C4").Value = 'Insurance (Life)'
$ws.Range("D4").Value = -0.0533
$ws.Range("G4").Value = -0.06782006920415225
$ws.Range("H4").Value = -0.06782006920415225
$ws.Range("I4").Value = -0.0342560553633218
$ws.Range("J4").Value = -0.0342560553633218
$ws.Range("K4").Value = -3.38
$ws.Range("L4").Value = -0.1169550173010381
$ws.Range("M4").Value = 1.66
$ws.Range("N4").Value = 0.04676056338028169
$ws.Range("O4").Value = -0.4911242603550296
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = 0.04676056338028169
$ws.Range("R4").Value = -0.4911242603550296
$ws.Range("S4").Value = 0.0
$ws.Range("T4").Value = 0.0
$ws.Range("U4").Value = 38.6
$ws.Range("V4").Value = 1.087323943661972
$ws.Range("W4").Value = -0.02857142857142857
$ws.Range("X4").Value = 0.08846408071440458
$ws.Range("Y4").Value = -0.1170355092858331
$ws.Range("Z4").Value = 0.366751269035533
$ws.Range("AA4").Value = -0.01256345177664975
$ws.Range("AB4").Value = 0.08846408071440458
$ws.Range("AC4").Value = -0.1010275324910543
$ws.Range("AD4").Value = 0.0
$ws.Range("AE4").Value = 0.0
$ws.Range("AF4").Value = 0.0
$ws.Range("AG4").Value = -38.6
$ws.Range("AH4").Value = 0.0
$ws.Range("AI4").Value = 0.0
$ws.Range("AJ4").Value = 12.4516129032258
$ws.Range("AK4").Value = -0.5065616797900263
$ws.Range("AL4").Value = 0.0
$ws.Range("AM4").Value = 0.0

# Row 5: Pragati Life Insurance Limited (DSE:PRAGATILIF)
$ws.Range("A5").Value = 'Bangladesh'
$ws.Range("B5").Value = 'Pragati Life Insurance Limited (DSE:PRAGATILIF)'
$ws.Range("C5").Value = 'Insurance (Life)'
$ws.Range("D5").Value = 0.08650000000000001
$ws.Range("E5").Value = -0.164
$ws.Range("G5").Value = 0.07005208333333333
$ws.Range("H5").Value = 0.07005208333333333
$ws.Range("I5").Value = 0.0609375
$ws.Range("J5").Value = 0.05163364955357143
$ws.Range("K5").Value = 1.9
$ws.Range("L5").Value = 0.04947916666666666
$ws.Range("M5").Value = 0.215
$ws.Range("N5").Value = 0.0134375
$ws.Range("O5").Value = 0.1131578947368421
$ws.Range("P5").Value = 0.215
$ws.Range("Q5").Value = 0.0134375
$ws.Range("R5").Value = 0.1131578947368421
$ws.Range("S5").Value = 0.0
$ws.Range("T5").Value = 0.0
$ws.Range("U5").Value = 26.8
$ws.Range("V5").Value = 1.675
$ws.Range("W5").Value = 1.043956043956044
$ws.Range("X5").Value = 0.08846408071440458
$ws.Range("Y5").Value = 0.9554919632416392
$ws.Range("Z5").Value = -1.513002364066194
$ws.Range("AA5").Value = -0.07812183383991894
$ws.Range("AB5").Value = 0.08846408071440458
$ws.Range("AC5").Value = -0.1665859145543235
$ws.Range("AD5").Value = 0.0
$ws.Range("AE5").Value = 0.0
$ws.Range("AF5").Value = 0.0
$ws.Range("AG5").Value = -26.8
$ws.Range("AH5").Value = 0.0
$ws.Range("AI5").Value = 0.0
$ws.Range("AJ5").Value = 2.481481481481481
$ws.Range("AK5").Value = 1.072428971588635
$ws.Range("AL5").Value = 0.0
$ws.Range("AM5").Value = 0.0
$ws.Range("AN5").Value = 0.0
$ws.Range("AP5").Value = -10.67729083665339

# Row 6: Meghna Life Insurance Company Limited (DSE:MEGHNALIFE)
$ws.Range("A6").Value = 'Bangladesh'
$ws.Range("B6").Value = 'Meghna Life Insurance Company Limited (DSE:MEGHNALIFE)'
$ws.Range("C6").Value = 'Insurance (Life)'
$ws.Range("D6").Value = -0.0106
$ws.Range("E6").Value = 0.0206
$ws.Range("G6").Value = 0.2191103789126853
$ws.Range("H6").Value = 0.2191103789126853
$ws.Range("I6").Value = 0.200988467874794
$ws.Range("J6").Value = 0.1859791677738215
$ws.Range("K6").Value = 11.5
$ws.Range("L6").Value = 0.1894563426688632
$ws.Range("M6").Value = 0.741
$ws.Range("N6").Value = 0.02975903614457831
$ws.Range("O6").Value = 0.06443478260869565
$ws.Range("P6").Value = 0.741
$ws.Range("Q6").Value = 0.02975903614457831
$ws.Range("R6").Value = 0.06443478260869565
$ws.Range("S6").Value = 0.0
$ws.Range("T6").Value = 0.0
$ws.Range("U6").Value = 97.8
$ws.Range("V6").Value = 3.927710843373494
$ws.Range("W6").Value = 2.896725440806045
$ws.Range("X6").Value = 0.08846408071440458
$ws.Range("Y6").Value = 2.80826136009164
$ws.Range("Z6").Value = -0.7014908124349937
$ws.Range("AA6").Value = -0.130462677497642
$ws.Range("AB6").Value = 0.08846408071440458
$ws.Range("AC6").Value = -0.2189267582120466
$ws.Range("AD6").Value = 0.0
$ws.Range("AE6").Value = 0.0
$ws.Range("AF6").Value = 0.0
$ws.Range("AG6").Value = -97.8
$ws.Range("AH6").Value = 0.0
$ws.Range("AI6").Value = 0.0
$ws.Range("AJ6").Value = 1.34156378600823
$ws.Range("AK6").Value = 1.042088438998402
$ws.Range("AL6").Value = 0.0
$ws.Range("AM6").Value = 0.0
$ws.Range("AN6").Value = 0.0
$ws.Range("AP6").Value = -7.761904761904762

# Row 7: National Life Insurance Company Limited (DSE:NATLIFEINS)
$ws.Range("A7").Value = 'Bangladesh'
$ws.Range("B7").Value = 'National Life Insurance Company Limited (DSE:NATLIFEINS)'
$ws.Range("C7").Value = 'Insurance (Life)'
$ws.Range("G7").Value = 0.2385697538100821
$ws.Range("H7").Value = 0.2385697538100821
$ws.Range("I7").Value = 0.2620164126611958
$ws.Range("J7").Value = 0.2347355275900644
$ws.Range("K7").Value = 39.1
$ws.Range("L7").Value = 0.2291910902696366
$ws.Range("M7").Value = 5.03
$ws.Range("N7").Value = 0.01542944785276074
$ws.Range("O7").Value = 0.1286445012787724
$ws.Range("P7").Value = 5.03
$ws.Range("Q7").Value = 0.01542944785276074
$ws.Range("R7").Value = 0.1286445012787724
$ws.Range("S7").Value = 0.0
$ws.Range("T7").Value = 0.0
$ws.Range("U7").Value = 210.3
$ws.Range("V7").Value = 0.6450920245398774
$ws.Range("W7").Value = 2.538961038961039
$ws.Range("X7").Value = 0.08846408071440458
$ws.Range("Y7").Value = 2.450496958246634
$ws.Range("Z7").Value = -0.8848547717842324
$ws.Range("AA7").Value = -0.2077068516953579
$ws.Range("AB7").Value = 0.08846408071440458
$ws.Range("AC7").Value = -0.2961709324097624
$ws.Range("AD7").Value = 0.0
$ws.Range("AE7").Value = 0.0
$ws.Range("AF7").Value = 0.0
$ws.Range("AG7").Value = -210.3
$ws.Range("AH7").Value = 0.0
$ws.Range("AI7").Value = 0.0
$ws.Range("AJ7").Value = -1.817631806395852
$ws.Range("AK7").Value = 1.082346886258363
$ws.Range("AL7").Value = 0.047
$ws.Range("AM7").Value = 0.047
$ws.Range("AN7").Value = 0.0
$ws.Range("AO7").Value = 951.0638297872341
$ws.Range("AP7").Value = -4.561822125813449
$ws.Range("AQ7").Value = 951.0638297872341
